$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: the two empty paragraphs right after "/create-tag: ..." become two
# new ListParagraph (numId=2) bullet items:
#   /delete-workspace: ...
#   /edit-workspace: ...
# ---------------------------------------------------------------------------

function Find-ParaIndexByText($doc, $needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

$idxCreateTag = Find-ParaIndexByText $d "/create-tag:"
$pCreateTag = $d.Paragraphs.Item($idxCreateTag)

# New paragraph 1: /delete-workspace...
$null = $pCreateTag.Range.InsertParagraphAfter()
$pNew1 = $d.Paragraphs.Item($idxCreateTag + 1)
$pNew1.Range.InsertAfter("/delete-workspace: only admin of the workspace can delete it (change status to ‘removed’)")

# New paragraph 2: /edit-workspace...
$null = $pNew1.Range.InsertParagraphAfter()
$pNew2 = $d.Paragraphs.Item($idxCreateTag + 2)
$pNew2.Range.InsertAfter("/edit-workspace: only admin can edit the workspace details.")

# The two original empty paragraphs got pushed down right after $pNew2; delete
# them now (delete the later one first so indices stay valid).
$pOldEmpty2 = $d.Paragraphs.Item($idxCreateTag + 4)
$pOldEmpty2.Range.Delete()
$pOldEmpty1 = $d.Paragraphs.Item($idxCreateTag + 3)
$pOldEmpty1.Range.Delete()

# ---------------------------------------------------------------------------
# Hunk 2: after "/mention-user-in-comment: ... workspace." insert four new
# ListParagraph (numId=3) bullet items:
#   /edit-project: ...
#   /delete-project: ... (two runs: text + ".")
#   /edit-task: ...
#   /delete-task: ... (three runs: text + "change status to 'removed'" + ").")
# ---------------------------------------------------------------------------

$idxMention = Find-ParaIndexByText $d "/mention-user-in-comment:"
$pMention = $d.Paragraphs.Item($idxMention)

# /edit-project
$null = $pMention.Range.InsertParagraphAfter()
$pEditProject = $d.Paragraphs.Item($idxMention + 1)
$pEditProject.Range.InsertAfter("/edit-project: user (admin or member) can edit project details and statuses.")

# /delete-project -- built as two paragraphs then joined so they stay as two
# separate runs (mirrors how the original document already splits runs).
$null = $pEditProject.Range.InsertParagraphAfter()
$pDeleteProjectA = $d.Paragraphs.Item($idxMention + 2)
$pDeleteProjectA.Range.InsertAfter("/delete-project: user (admin or member) can delete a project (change status to ‘removed’)")

$null = $pDeleteProjectA.Range.InsertParagraphAfter()
$pDeleteProjectB = $d.Paragraphs.Item($idxMention + 3)
$pDeleteProjectB.Range.InsertAfter(".")

$joinMark = $d.Range($pDeleteProjectA.Range.End - 1, $pDeleteProjectA.Range.End)
$joinMark.Delete()
# $pDeleteProjectA now holds the full "/delete-project...'removed')." text as two runs.

# /edit-task
$null = $pDeleteProjectA.Range.InsertParagraphAfter()
$pEditTask = $d.Paragraphs.Item($idxMention + 3)
$pEditTask.Range.InsertAfter("/edit-task: user (admin or member) can edit task details and tags.")

# /delete-task -- built as three paragraphs then joined back-to-front so they
# stay as three separate runs.
$null = $pEditTask.Range.InsertParagraphAfter()
$pDeleteTaskA = $d.Paragraphs.Item($idxMention + 4)
$pDeleteTaskA.Range.InsertAfter("/delete-task: user (admin or member) can delete a task (")

$null = $pDeleteTaskA.Range.InsertParagraphAfter()
$pDeleteTaskB = $d.Paragraphs.Item($idxMention + 5)
$pDeleteTaskB.Range.InsertAfter("change status to ‘removed’")

$null = $pDeleteTaskB.Range.InsertParagraphAfter()
$pDeleteTaskC = $d.Paragraphs.Item($idxMention + 6)
$pDeleteTaskC.Range.InsertAfter(").")

$joinMark2 = $d.Range($pDeleteTaskB.Range.End - 1, $pDeleteTaskB.Range.End)
$joinMark2.Delete()
$joinMark1 = $d.Range($pDeleteTaskA.Range.End - 1, $pDeleteTaskA.Range.End)
$joinMark1.Delete()

Write-Host "Done."
